# Commit: "Update gh-pages to output generated at 456a3b4"
# Refreshes the scraped "想去人数" (want-to-go count) column F — and, on two
# rows, the "最低票价" (min ticket price) column G — across all four sheets
# of the 上海-漫展信息 workbook with newly-scraped numbers.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 159
$ws.Cells.Item(3, 6).Value = 975
$ws.Cells.Item(4, 6).Value = 610
$ws.Cells.Item(5, 6).Value = 3044
$ws.Cells.Item(6, 6).Value = 810
$ws.Cells.Item(7, 6).Value = 601
$ws.Cells.Item(8, 6).Value = 606
$ws.Cells.Item(9, 6).Value = 456
$ws.Cells.Item(10, 6).Value = 661
$ws.Cells.Item(12, 6).Value = 574
$ws.Cells.Item(14, 6).Value = 2165
$ws.Cells.Item(16, 6).Value = 751
$ws.Cells.Item(17, 6).Value = 41
$ws.Cells.Item(19, 6).Value = 2688
$ws.Cells.Item(19, 7).Value = 39.9
$ws.Cells.Item(20, 6).Value = 14
$ws.Cells.Item(23, 6).Value = 536
$ws.Cells.Item(25, 6).Value = 676
$ws.Cells.Item(26, 6).Value = 676
$ws.Cells.Item(30, 6).Value = 24
$ws.Cells.Item(31, 6).Value = 559
$ws.Cells.Item(33, 6).Value = 230
$ws.Cells.Item(34, 6).Value = 128
$ws.Cells.Item(35, 6).Value = 913
$ws.Cells.Item(36, 6).Value = 4712
$ws.Cells.Item(37, 6).Value = 270
$ws.Cells.Item(38, 6).Value = 44
$ws.Cells.Item(39, 6).Value = 13

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 6
$ws.Cells.Item(9, 6).Value = 364
$ws.Cells.Item(11, 6).Value = 2
$ws.Cells.Item(14, 6).Value = 172
$ws.Cells.Item(16, 6).Value = 278
$ws.Cells.Item(20, 6).Value = 63
$ws.Cells.Item(23, 6).Value = 277
$ws.Cells.Item(24, 6).Value = 25
$ws.Cells.Item(25, 6).Value = 309
$ws.Cells.Item(27, 6).Value = 257
$ws.Cells.Item(31, 6).Value = 33
$ws.Cells.Item(33, 6).Value = 5
$ws.Cells.Item(37, 6).Value = 601

# 本地生活 (Local Life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 1470
$ws.Cells.Item(5, 6).Value = 579
$ws.Cells.Item(6, 6).Value = 281
$ws.Cells.Item(7, 6).Value = 275

# 全部类型 (All Types) - combined view, same rows as above 3 sheets concatenated
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1470
$ws.Cells.Item(4, 6).Value = 579
$ws.Cells.Item(5, 6).Value = 159
$ws.Cells.Item(6, 6).Value = 281
$ws.Cells.Item(7, 6).Value = 975
$ws.Cells.Item(8, 6).Value = 610
$ws.Cells.Item(9, 6).Value = 3044
$ws.Cells.Item(10, 6).Value = 810
$ws.Cells.Item(11, 6).Value = 601
$ws.Cells.Item(12, 6).Value = 606
$ws.Cells.Item(13, 6).Value = 456
$ws.Cells.Item(14, 6).Value = 661
$ws.Cells.Item(15, 6).Value = 6
$ws.Cells.Item(17, 6).Value = 574
$ws.Cells.Item(18, 6).Value = 364
$ws.Cells.Item(20, 6).Value = 2
$ws.Cells.Item(21, 6).Value = 2165
$ws.Cells.Item(23, 6).Value = 751
$ws.Cells.Item(24, 6).Value = 41
$ws.Cells.Item(27, 6).Value = 2688
$ws.Cells.Item(27, 7).Value = 39.9
$ws.Cells.Item(28, 6).Value = 14
$ws.Cells.Item(29, 6).Value = 278
$ws.Cells.Item(32, 6).Value = 536
$ws.Cells.Item(33, 6).Value = 275
$ws.Cells.Item(35, 6).Value = 676
$ws.Cells.Item(36, 6).Value = 676
$ws.Cells.Item(38, 6).Value = 277
$ws.Cells.Item(39, 6).Value = 24
$ws.Cells.Item(40, 6).Value = 559
$ws.Cells.Item(42, 6).Value = 309
$ws.Cells.Item(43, 6).Value = 230
$ws.Cells.Item(45, 6).Value = 128
$ws.Cells.Item(46, 6).Value = 913
$ws.Cells.Item(47, 6).Value = 4712
$ws.Cells.Item(48, 6).Value = 270
$ws.Cells.Item(49, 6).Value = 44
$ws.Cells.Item(50, 6).Value = 601
$ws.Cells.Item(51, 6).Value = 601
